$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number and must be forced to text
# to preserve the original inline-string formatting (e.g. "1.00", "215.26").
$textForceCells = @(
    "D5","D9","D11","D16","D18","D19","D22","D23","D25","D26","D28","D29","D32","D36","D37","D40","D42","D46","D47","D49","D50","D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.033.14'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.676.30'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '215.26'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D9').Value = '21.42'
$ws.Range('E9').Value = '  +5.49%  '
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '1.912.18'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.675.14'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '66.33'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '27.030.77'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '8.20'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('D19').Value = '236.00'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = '4.47'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').Value = '9.25'
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('E24').Value = '  -2.46%  '
$ws.Range('D25').Value = '148.09'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '7.28'
$ws.Range('E26').Value = '  +2.13%  '
$ws.Range('E27').Value = '  +3.68%  '
$ws.Range('D28').Value = '0.113'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').Value = '3.38'
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').Value = '1.542.84'
$ws.Range('E33').Value = '  +6.86%  '
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('E35').Value = '  +5.22%  '
$ws.Range('D36').Value = '2.39'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Value = '0.589'
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('D40').Value = '1.05'
$ws.Range('E40').Value = '  +5.12%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = '67.90'
$ws.Range('E42').Value = '  +2.95%  '
$ws.Range('E43').Value = '  -3.42%  '
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('D45').Value = '1.819.70'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('D46').Value = '0.780'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('D47').Value = '90.36'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').Value = '0.104'
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('D50').Value = '8.04'
$ws.Range('E50').Value = '  +6.65%  '
$ws.Range('D51').Value = '0.0509'
$ws.Range('E51').Value = '  +0.38%  '
